# Update "How to talk to SUMO through python API.pptx":
#  - Slide 3 ("Please download the simulation files ...")'s TextBox 3:
#      * 1st bullet: replace the SharePoint link text with the new GitHub
#        "Data" link and bump its font size from 12pt to 16pt.
#      * 2nd bullet: merge the two runs ("Script can be found here" and
#        ": https://.../Code") into a single run with one consistent format.
#      * The shape uses shrink-to-fit autosizing, so its height is
#        recalculated to match the new (shorter) wrapped text.

$p = $ppt.ActivePresentation

# Locate the slide/shape by name instead of a hard-coded index so the
# script is resilient to any slide-order differences.
$targetSlide = $null
$targetShape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    for ($j = 1; $j -le $sl.Shapes.Count; $j++) {
        $sh = $sl.Shapes.Item($j)
        if ($sh.Name -eq "TextBox 3") {
            $targetSlide = $sl
            $targetShape = $sh
        }
    }
}

$shape = $targetShape
$tr = $shape.TextFrame.TextRange

# --- Bullet 1: "Simulation files can be found here: <link>" -----------
# Only the link run (2nd run of paragraph 1) changes: new URL + 16pt size.
$para1 = $tr.Paragraphs(1)
$p1Text = $para1.Text
$urlIdx = $p1Text.IndexOf("https://")
$urlStart = $para1.Start + $urlIdx
$urlLen = $para1.Length - $urlIdx
$urlRun = $tr.Characters($urlStart, $urlLen)
# Set formatting before text so the Characters() range (fixed start/length)
# still lines up with the old (longer) run -- changing Text first would
# shrink the paragraph and make the stale range bleed into later text.
$urlRun.Font.Size = 16
$urlRun.Text = "https://github.com/ivsg-psu/TrafficSimulators_GettingStartedWithDifferrentSimulators_GettingStartedWithSUMO/tree/main/Data"

# --- Bullet 2: "Script can be found here: <link>" ----------------------
# Originally two runs ("Script can be found here" / ": <link>"); collapse
# to a single run using the first run's formatting.
$para2 = $tr.Paragraphs(2)
$oldLen2 = $para2.Length
$label = "Script can be found here"
$run1 = $tr.Characters($para2.Start, $label.Length)
$newText2 = "Script can be found here: https://github.com/ivsg-psu/TrafficSimulators_GettingStartedWithDifferrentSimulators_GettingStartedWithSUMO/tree/main/Code"
$run1.Text = $newText2
# Remove the leftover tail of the old 2nd run (now sitting right after the
# freshly written text).
$remainderStart = $para2.Start + $newText2.Length
$remainderLen = $oldLen2 - $label.Length
$remainder = $tr.Characters($remainderStart, $remainderLen)
$remainder.Text = ""

# --- Shape autosize: shrink-to-fit box now wraps less text -------------
$shape.Height = 1723549 / 914400 * 72
